# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns
# with freshly scraped values, as produced by the scheduled GitHub
# Actions job that refreshes cryptos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ Price = <new D value, or $null if the price cell is unchanged>; Volume = <new E value> }
$updates = @{
    2 = @{ Price = "26.421.68"; Volume = "  -1.67%  " }
    3 = @{ Price = "1.841.18"; Volume = "  -2.07%  " }
    4 = @{ Price = "1.000"; Volume = "  -0.12%  " }
    5 = @{ Price = "259.88"; Volume = "  -6.67%  " }
    6 = @{ Price = "1.000"; Volume = "  -0.08%  " }
    7 = @{ Price = "0.5200"; Volume = "  -2.14%  " }
    8 = @{ Price = "0.3261"; Volume = "  -5.59%  " }
    9 = @{ Price = "0.06771"; Volume = "  -2.83%  " }
    10 = @{ Price = "18.56"; Volume = "  -7.71%  " }
    11 = @{ Price = "0.7649"; Volume = "  -5.15%  " }
    12 = @{ Price = "0.07701"; Volume = "  -0.43%  " }
    13 = @{ Price = "1.826.84"; Volume = "  -2.88%  " }
    14 = @{ Price = "88.10"; Volume = "  -2.67%  " }
    15 = @{ Price = "5.022"; Volume = "  -3.07%  " }
    16 = @{ Price = $null; Volume = "  -0.04%  " }
    17 = @{ Price = "13.87"; Volume = "  -4.80%  " }
    18 = @{ Price = $null; Volume = "  +0.07%  " }
    19 = @{ Price = "0.000007950"; Volume = "  -1.07%  " }
    20 = @{ Price = "26.442.75"; Volume = "  -1.81%  " }
    21 = @{ Price = "2.073.84"; Volume = "  -2.32%  " }
    22 = @{ Price = "4.557"; Volume = "  -4.19%  " }
    23 = @{ Price = "9.460"; Volume = "  -5.89%  " }
    24 = @{ Price = "5.948"; Volume = "  -4.28%  " }
    25 = @{ Price = "144.43"; Volume = "  -1.80%  " }
    26 = @{ Price = "2.213"; Volume = "  -6.49%  " }
    27 = @{ Price = "1.642"; Volume = "  -1.21%  " }
    28 = @{ Price = "16.96"; Volume = "  -2.29%  " }
    29 = @{ Price = "111.14"; Volume = "  -2.17%  " }
    30 = @{ Price = "4.158"; Volume = "  -4.63%  " }
    31 = @{ Price = "4.118"; Volume = "  -4.76%  " }
    32 = @{ Price = "0.08711"; Volume = "  -2.07%  " }
    33 = @{ Price = "0.04772"; Volume = "  -3.33%  " }
    34 = @{ Price = "1.123"; Volume = "  -4.36%  " }
    35 = @{ Price = "2.849"; Volume = "  -1.68%  " }
    36 = @{ Price = "0.7011"; Volume = "  -4.37%  " }
    37 = @{ Price = "3.053"; Volume = "  -7.38%  " }
    38 = @{ Price = "0.01749"; Volume = "  -5.62%  " }
    39 = @{ Price = "2.180"; Volume = "  -8.17%  " }
    40 = @{ Price = "0.4810"; Volume = "  -6.52%  " }
    41 = @{ Price = "110.93"; Volume = "  -4.13%  " }
    42 = @{ Price = "0.8927"; Volume = "  -6.78%  " }
    43 = @{ Price = "6.063"; Volume = "  -2.04%  " }
    44 = @{ Price = "1.000"; Volume = "  -0.02%  " }
    45 = @{ Price = "7.653"; Volume = "  -5.90%  " }
    46 = @{ Price = "0.05867"; Volume = "  -1.40%  " }
    47 = @{ Price = "0.4104"; Volume = "  -8.33%  " }
    48 = @{ Price = "8.954"; Volume = "  -4.70%  " }
    49 = @{ Price = $null; Volume = "  -3.45%  " }
    50 = @{ Price = "0.1215"; Volume = "  -9.55%  " }
    51 = @{ Price = "0.8861"; Volume = "  -0.01%  " }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($null -ne $u.Price) {
        # Force the Price cell to remain plain text so that values such as
        # "1.000" or "0.5200" keep their exact trailing zeros / grouping dots
        # instead of being re-interpreted as numbers by Excel.
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    $ws.Cells.Item($row, 5).Value = $u.Volume
}
